# Update the Password value for the last "NON_SSO" / QA row on LoginData
# from the old rotated credential to the new one, and move the active
# selection cursor from A5 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# F6 holds the password string - rotate it to the new value.
$ws.Range("F6").Value = "202512KFLog!"

# Update the saved selection/active cell on the LoginData sheet.
$ws.Activate()
$ws.Range("A4").Select()
